# Auto-generated edit script: updates the cryptos price table (rows 2-51)
# to match the latest GitHub Actions scrape. Values in column D that look
# like plain numbers are written with a leading apostrophe so Excel keeps
# them as text (matching the original inline-string cell type) instead of
# silently converting them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.669.99"
$ws.Range("E2").Value = "  +1.83%  "

$ws.Range("D3").Value = "1.656.76"
$ws.Range("E3").Value = "  +3.00%  "

$ws.Range("D4").Value = "'0.9973"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'306.16"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("D6").Value = "'0.9973"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "'0.3778"

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3673"
$ws.Range("E8").Value = "  +0.88%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'52.16"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("D10").Value = "'1.273"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").Value = "'0.08178"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("D12").Value = "'0.9972"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").Value = "'23.19"
$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("D14").Value = "'6.719"
$ws.Range("E14").Value = "  +1.84%  "

$ws.Range("D15").Value = "'0.00001278"
$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("D16").Value = "'7.406"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "1.663.25"
$ws.Range("E17").Value = "  +3.49%  "

$ws.Range("D18").Value = "'95.60"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").Value = "'0.06920"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").Value = "'18.49"

$ws.Range("D21").Value = "'6.615"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").Value = "'0.9978"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").Value = "23.672.79"
$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("D24").Value = "'12.96"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").Value = "'3.158"
$ws.Range("E25").Value = "  +3.49%  "

$ws.Range("D26").Value = "'2.417"
$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("D27").Value = "'21.47"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").Value = "'151.22"
$ws.Range("E28").Value = "  +1.01%  "

$ws.Range("D29").Value = "'5.310"
$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").Value = "'137.30"
$ws.Range("E30").Value = "  +1.21%  "

$ws.Range("D31").Value = "'2.335"
$ws.Range("E31").Value = "  -1.91%  "

$ws.Range("D32").Value = "1.845.70"
$ws.Range("E32").Value = "  +3.95%  "

$ws.Range("D33").Value = "'6.933"
$ws.Range("E33").Value = "  +2.41%  "

$ws.Range("D34").Value = "'11.07"
$ws.Range("E34").Value = "  +6.66%  "

$ws.Range("D35").Value = "'0.9802"
$ws.Range("E35").Value = "  +1.28%  "

$ws.Range("D36").Value = "'0.02887"
$ws.Range("E36").Value = "  +4.70%  "

$ws.Range("D37").Value = "'6.400"
$ws.Range("E37").Value = "  +4.20%  "

$ws.Range("D38").Value = "'0.2578"
$ws.Range("E38").Value = "  +2.21%  "

$ws.Range("D39").Value = "'0.07414"
$ws.Range("E39").Value = "  -1.24%  "

$ws.Range("D40").Value = "'0.08912"
$ws.Range("E40").Value = "  +1.13%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.7247"
$ws.Range("E41").Value = "  +1.79%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.382"
$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("D43").Value = "'16.46"
$ws.Range("E43").Value = "  +5.13%  "

$ws.Range("D44").Value = "'12.68"
$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("D45").Value = "'0.6673"
$ws.Range("E45").Value = "  +1.84%  "

$ws.Range("D46").Value = "'2.380"
$ws.Range("E46").Value = "  +1.89%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'0.9972"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'4.030"
$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("D49").Value = "'0.08052"
$ws.Range("E49").Value = "  +1.22%  "

$ws.Range("D50").Value = "'1.232"
$ws.Range("E50").Value = "  +1.86%  "

$ws.Range("D51").Value = "'128.99"
$ws.Range("E51").Value = "  -3.02%  "
